$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 11643
$ws.Range("I34").Value = 4734.857
$ws.Range("K34").Value = 4734.857
$ws.Range("M34").Value = -4531.857
$ws.Range("H36").Value = 11643
$ws.Range("I36").Value = 4734.857
$ws.Range("K36").Value = 4734.857
$ws.Range("M36").Value = -4019.857
$ws.Range("H109").Value = 16666.666
$ws.Range("J109").Value = 16666.666
$ws.Range("L109").Value = 16666.666
$ws.Range("N109").Value = -19440.666
$ws.Range("H113").Value = 103588.4
$ws.Range("I113").Value = 114664.89
$ws.Range("K113").Value = 114664.89
$ws.Range("M113").Value = -111410.89
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 125379.875
$ws.Range("I5").Value = 167039.83
$ws.Range("K5").Value = 167039.83
$ws.Range("M5").Value = -166927.83
$ws.Range("H46").Value = 6712.75
$ws.Range("J46").Value = 6712.75
$ws.Range("L46").Value = 6712.75
$ws.Range("N46").Value = -7350.75
$ws.Range("H74").Value = 3018.2686
$ws.Range("I74").Value = 909.95746
$ws.Range("J74").Value = 7972.8
$ws.Range("K74").Value = 909.95746
$ws.Range("L74").Value = 7972.8
$ws.Range("M74").Value = -35.95745999999997
$ws.Range("N74").Value = -9720.799999999999
$ws.Range("H77").Value = 3018.2686
$ws.Range("I77").Value = 909.95746
$ws.Range("J77").Value = 7972.8
$ws.Range("K77").Value = 4549.7873
$ws.Range("L77").Value = 39864
$ws.Range("M77").Value = -181.7873
$ws.Range("N77").Value = -48600
$ws.Range("H97").Value = 22229314
$ws.Range("I97").Value = 22229314
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 22229314
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -22228818
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 2842.1
$ws.Range("I122").Value = 2865.125
$ws.Range("K122").Value = 8595.375
$ws.Range("M122").Value = -6145.375
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 125379.875
$ws.Range("I4").Value = 167039.83
$ws.Range("K4").Value = 167039.83
$ws.Range("M4").Value = -166924.83
$ws.Range("H134").Value = 20002192
$ws.Range("I134").Value = 25642766
$ws.Range("J134").Value = 3794.5454
$ws.Range("K134").Value = 76928298
$ws.Range("L134").Value = 11383.6362
$ws.Range("M134").Value = -76925763
$ws.Range("N134").Value = -16453.6362
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2405.3076
$ws.Range("I5").Value = 236.14285
$ws.Range("J5").Value = 4936
$ws.Range("K5").Value = 236.14285
$ws.Range("L5").Value = 4936
$ws.Range("M5").Value = -124.14285
$ws.Range("N5").Value = -5160
$ws.Range("H41").Value = 14509.6
$ws.Range("I41").Value = 5806
$ws.Range("J41").Value = 27565
$ws.Range("K41").Value = 5806
$ws.Range("L41").Value = 27565
$ws.Range("M41").Value = -5378
$ws.Range("N41").Value = -28421
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 768.93335
$ws.Range("I8").Value = 768.93335
$ws.Range("K8").Value = 2306.80005
$ws.Range("M8").Value = -2167.80005
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H75").Value = 29900
$ws.Range("J75").Value = 29900
$ws.Range("L75").Value = 29900
$ws.Range("N75").Value = -31648
$ws.Range("H78").Value = 29900
$ws.Range("J78").Value = 29900
$ws.Range("L78").Value = 89700
$ws.Range("N78").Value = -98436
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2992.3809
$ws.Range("I7").Value = 2050
$ws.Range("J7").Value = 3369.3333
$ws.Range("K7").Value = 2050
$ws.Range("L7").Value = 3369.3333
$ws.Range("M7").Value = -1938
$ws.Range("N7").Value = -3593.3333
$ws.Range("H17").Value = 1999.8
$ws.Range("J17").Value = 2374.75
$ws.Range("L17").Value = 2374.75
$ws.Range("N17").Value = -2714.75
$ws.Range("H26").Value = 19333.334
$ws.Range("I26").Value = 26500
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 26500
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = -26205
$ws.Range("N26").Value = -5590
$ws.Range("H31").Value = 1086
$ws.Range("I31").Value = 676.6667
$ws.Range("K31").Value = 676.6667
$ws.Range("M31").Value = -428.6667
$ws.Range("H53").Value = 5720
$ws.Range("I53").Value = 3000
$ws.Range("K53").Value = 3000
$ws.Range("M53").Value = -2482
$ws.Range("H122").Value = 3833.3333
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3833.3333
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11499.9999
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -16399.9999
$ws.Range("H126").Value = 2992.3809
$ws.Range("I126").Value = 2050
$ws.Range("J126").Value = 3369.3333
$ws.Range("K126").Value = 6150
$ws.Range("L126").Value = 10107.9999
$ws.Range("M126").Value = -3680
$ws.Range("N126").Value = -15047.9999
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 10664.833
$ws.Range("I9").Value = 26000
$ws.Range("J9").Value = 2997.25
$ws.Range("K9").Value = 26000
$ws.Range("L9").Value = 2997.25
$ws.Range("M9").Value = -25860
$ws.Range("N9").Value = -3277.25
$ws.Range("H30").Value = 7500
$ws.Range("J30").Value = 7500
$ws.Range("L30").Value = 7500
$ws.Range("N30").Value = -7714
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H55").Value = 1156.7142
$ws.Range("I55").Value = 661
$ws.Range("J55").Value = 1817.6666
$ws.Range("K55").Value = 661
$ws.Range("L55").Value = 1817.6666
$ws.Range("M55").Value = -384
$ws.Range("N55").Value = -2371.6666
$ws.Range("H92").Value = 31950
$ws.Range("J92").Value = 31950
$ws.Range("L92").Value = 31950
$ws.Range("N92").Value = -36942
$ws.Range("H122").Value = 102400.4
$ws.Range("I122").Value = 251501
$ws.Range("K122").Value = 754503
$ws.Range("M122").Value = -752053
$ws.Range("H136").Value = 15163.043
$ws.Range("I136").Value = 16651.555
$ws.Range("K136").Value = 49954.665
$ws.Range("M136").Value = -47404.665
